# Applies the cryptos list price/volume update described in the commit
# "Updated cryptos list on Tue Aug  1 13:07:27 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# are forced to keep a Text format so the value round-trips as a string,
# matching the inline-string cells produced by the source export.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D13", "D14", "D15", "D16", "D17", "D19", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.852.98'
$ws.Range("E2").Value = '  -2.04%  '
$ws.Range("D3").Value = '1.830.79'
$ws.Range("E3").Value = '  -2.18%  '
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '245.13'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").Value = '0.6899'
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.07702'
$ws.Range("E8").Value = '  -3.20%  '
$ws.Range("D9").Value = '0.3051'
$ws.Range("E9").Value = '  -3.03%  '
$ws.Range("E10").Value = '  -4.98%  '
$ws.Range("D11").Value = '0.07804'
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").Value = '1.829.56'
$ws.Range("E12").Value = '  -2.81%  '
$ws.Range("D13").Value = '5.088'
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").Value = '90.69'
$ws.Range("E14").Value = '  -3.34%  '
$ws.Range("D15").Value = '0.6815'
$ws.Range("D16").Value = '6.442'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = '0.000008299'
$ws.Range("E17").Value = '  -3.43%  '
$ws.Range("D18").Value = '28.846.25'
$ws.Range("E18").Value = '  -2.23%  '
$ws.Range("D19").Value = '241.66'
$ws.Range("E19").Value = '  -4.21%  '
$ws.Range("D20").Value = '2.076.45'
$ws.Range("E20").Value = '  -3.34%  '
$ws.Range("D21").Value = '12.69'
$ws.Range("E21").Value = '  -3.27%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '7.464'
$ws.Range("E23").Value = '  -2.29%  '
$ws.Range("D24").Value = '0.9999'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = '0.1483'
$ws.Range("E25").Value = '  -4.06%  '
$ws.Range("D26").Value = '158.70'
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("D27").Value = '8.789'
$ws.Range("E27").Value = '  -2.50%  '
$ws.Range("E28").Value = '  -2.94%  '
$ws.Range("D29").Value = '1.549'
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("D30").Value = '4.216'
$ws.Range("E30").Value = '  -2.28%  '
$ws.Range("D31").Value = '4.148'
$ws.Range("E31").Value = '  -2.85%  '
$ws.Range("D32").Value = '1.183'
$ws.Range("D33").Value = '0.05111'
$ws.Range("E33").Value = '  -3.48%  '
$ws.Range("D34").Value = '0.7725'
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").Value = '1.849'
$ws.Range("E35").Value = '  -2.55%  '
$ws.Range("D36").Value = '1.140'
$ws.Range("E36").Value = '  -3.86%  '
$ws.Range("D37").Value = '2.687'
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("D38").Value = '0.01847'
$ws.Range("E38").Value = '  -1.88%  '
$ws.Range("D39").Value = '1.225.19'
$ws.Range("E39").Value = '  -3.86%  '
$ws.Range("D40").Value = '2.693'
$ws.Range("E40").Value = '  -2.60%  '
$ws.Range("D41").Value = '0.9539'
$ws.Range("E41").Value = '  +5.68%  '
$ws.Range("D42").Value = '107.89'
$ws.Range("E42").Value = '  -1.66%  '
$ws.Range("D43").Value = '5.830'
$ws.Range("E43").Value = '  -2.87%  '
$ws.Range("D44").Value = '0.9993'
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = '9.638'
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("E46").Value = '  -3.46%  '
$ws.Range("D47").Value = '1.975.22'
$ws.Range("E47").Value = '  -3.48%  '
$ws.Range("D48").Value = '0.5154'
$ws.Range("E48").Value = '  -0.35%  '
$ws.Range("D49").Value = '64.39'
$ws.Range("E49").Value = '  -9.39%  '
$ws.Range("E50").Value = '  -3.30%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '6.914'
$ws.Range("E51").Value = '  -2.32%  '
